$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: finish off the existing (previously blank) label cell and add a value next to it
$ws.Range("A13").Value = "id "
$ws.Range("A13").Font.Bold = $true
$ws.Range("B13").Value = 2

# Row 14: full_navn / Administrator
$ws.Range("A14").Value = "full_navn"
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").Value = "Administrator"

# Row 15: brukernavn / admin
$ws.Range("A15").Value = "brukernavn"
$ws.Range("A15").Font.Bold = $true
$ws.Range("B15").Value = "admin "

# Row 16: passord / admin05
$ws.Range("A16").Value = "passord"
$ws.Range("A16").Font.Bold = $true
$ws.Range("B16").Value = "admin05"
